$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the paragraph that ends with "...Dining room." and the
# "Backstory." title paragraph that immediately follows it.
$housePara = $null
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Dining room.*") {
        $housePara = $i
    }
    if ($t -like "Backstory.*") {
        $titlePara = $i
        break
    }
}

# 1. Insert a brand new, completely blank paragraph right after the
#    "House contains..." paragraph (and before the "Backstory." title
#    paragraph). Using raw OOXML injection keeps it a bare <w:p/> with
#    no inherited style/numbering/run.
$house = $d.Paragraphs($housePara)
$insertionPoint = $d.Range($house.Range.End, $house.Range.End)
[void]$insertionPoint.InsertXML("<w:p $wNs/>")

# The title paragraph shifted down by one because of the insertion above.
$titlePara = $titlePara + 1
$title = $d.Paragraphs($titlePara)

# 2. Rename "Backstory." to "Character Profile" (keeps the Title style).
$title.Range.Text = "Character Profile"
$title = $d.Paragraphs($titlePara)

# 3. Insert the new Character Profile paragraphs right after the title,
#    one at a time, resetting each to the plain "Normal" style so no
#    Title formatting carries over.
$title.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($titlePara + 1)
$p1.Style = "Normal"
$p1.Range.Text = "Graham Sterling"

$p1 = $d.Paragraphs($titlePara + 1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($titlePara + 2)
$p2.Style = "Normal"
$p2.Range.Text = "76-year-old man living alone after his wife" + [char]0x2019 + "s death."

$p2 = $d.Paragraphs($titlePara + 2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($titlePara + 3)
$p3.Style = "Normal"
$p3.Range.Text = "Severe dementia"

# 4. Remove the trailing blank paragraph that used to sit just before
#    the section properties (end of document). The very last paragraph
#    mark of the body can't be deleted on its own, so fold it into the
#    previous paragraph's mark instead.
$n = $d.Paragraphs.Count
$lastText = $d.Paragraphs($n).Range.Text
if ($lastText -eq "`r" -or $lastText -eq "") {
    $prevEnd = $d.Paragraphs($n - 1).Range.End
    $lastEnd = $d.Paragraphs($n).Range.End
    $d.Range($prevEnd - 1, $lastEnd).Delete()
}
